# Repull data, push all data, mean calculation
# Update the dSF (F column) values for the affected rows to match the
# repulled / recalculated data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    5  = 0
    6  = 2
    10 = 0
    11 = 1
    21 = -3
    22 = -3
    23 = 0
    30 = -7
    32 = 0
    33 = -2
    37 = 3
    39 = -3
    50 = -4
    51 = -4
    53 = -4
}

foreach ($row in $updates.Keys) {
    $ws.Range("F$row").Value = $updates[$row]
}
